$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each touched cell is forced to the "Text" number format before its value is
# written, so Excel keeps it stored as a literal string (matching the original
# inline-string cell) instead of silently re-interpreting a numeric-looking or
# percent-looking string as a number. (NumberFormat is set per-cell rather than
# via a multi-area Union range, since only the first area of a Union reliably
# picks up the format in this COM host.)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "314.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.92%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "20"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.95"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-7.38%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "20"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.132"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.93%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "20"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07863"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.48%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "20"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.349"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.18%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "20"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.659"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-13.68%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "20"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9259"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.49%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "20"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1070"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.30%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "20"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1799"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.88%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "20"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09071"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.59%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "20"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04470"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.40%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "20"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.249"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-16.96%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "20"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1062"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.14%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "20"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001271"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.31%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "20"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005905"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.50%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "20"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.42%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "20"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.77%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.06%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "20"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1386"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.59%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "20"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2646"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.56%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "20"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04185"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.51%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "20"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001251"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.27%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-6.10%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "20"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001228"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.81%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "20"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003015"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.86%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "20"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "20"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "20"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "20"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "20"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "20"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "20"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "20"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "20"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "20"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "20"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "20"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02452"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-9.26%"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "20"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05332"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.96%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "20"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.008001"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.06%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "20"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.42%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "20"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007664"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.08%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "20"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001891"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-10.95%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "20"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008248"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.76%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "20"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3124"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-10.62%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "20"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006813"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.84%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "20"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.90%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "20"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-2.07%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "20"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004149"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "17.14%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "20"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002125"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.90%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "20"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002024"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.90%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "20"
